# -----------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (and before
#    "总计") holding the per-fund holdings detail table for the new
#    quarter (same layout as the other quarterly sheets).
# 2. Prepend a corresponding summary row to the "总计" (totals) sheet and
#    renumber the existing index column.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# Step 1: create the new "2022-Q1" worksheet, positioned right after
# the existing "2021-Q4" sheet (i.e. right before "总计").
# -------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$ws = $wb.Worksheets.Add($null, $q4)
$ws.Name = "2022-Q1"

# Clone the header-row / index-column formatting from the "2021-Q4"
# sheet (identical table layout) so the new sheet's styling matches
# the other quarter sheets exactly (bold + bordered + centred cells).
$q4.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$q4.Range("A2").Copy()
$ws.Range("A2:A18").PasteSpecial(-4122)

# Header row
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Data rows 2-18. Columns B, D, E, F, G are numeric-looking text (fund
# codes / percentages) and must stay text, so they are entered with a
# leading apostrophe - exactly like the other quarter sheets store
# them (t="inlineStr"/t="s"), instead of being auto-coerced into
# numbers. Column H (rank) is a genuine number.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "'011708"
$ws.Range("C2").Value = "中欧嘉益一年持有期混合型证券投资基金A"
$ws.Range("D2").Value = "'8.11"
$ws.Range("E2").Value = "'80.96"
$ws.Range("F2").Value = "'5.50"
$ws.Range("G2").Value = "'0.4460"
$ws.Range("H2").Value = 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "'005421"
$ws.Range("C3").Value = "中欧嘉泽灵活配置混合"
$ws.Range("D3").Value = "'8.87"
$ws.Range("E3").Value = "'86.45"
$ws.Range("F3").Value = "'4.98"
$ws.Range("G3").Value = "'0.4417"
$ws.Range("H3").Value = 3
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "'009230"
$ws.Range("C4").Value = "鹏华安和混合A"
$ws.Range("D4").Value = "'14.02"
$ws.Range("E4").Value = "'34.45"
$ws.Range("F4").Value = "'1.67"
$ws.Range("G4").Value = "'0.2341"
$ws.Range("H4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "'009667"
$ws.Range("C5").Value = "鹏华安庆混合A"
$ws.Range("D5").Value = "'11.22"
$ws.Range("E5").Value = "'38.92"
$ws.Range("F5").Value = "'1.81"
$ws.Range("G5").Value = "'0.2031"
$ws.Range("H5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "'011709"
$ws.Range("C6").Value = "中欧嘉益一年持有期混合型证券投资基金C"
$ws.Range("D6").Value = "'2.43"
$ws.Range("E6").Value = "'80.96"
$ws.Range("F6").Value = "'5.50"
$ws.Range("G6").Value = "'0.1336"
$ws.Range("H6").Value = 3
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "'550001"
$ws.Range("C7").Value = "信诚四季红混合"
$ws.Range("D7").Value = "'5.00"
$ws.Range("E7").Value = "'72.84"
$ws.Range("F7").Value = "'2.58"
$ws.Range("G7").Value = "'0.1290"
$ws.Range("H7").Value = 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "'005416"
$ws.Range("C8").Value = "鹏华尊惠18个月定期开放混合A"
$ws.Range("D8").Value = "'7.95"
$ws.Range("E8").Value = "'37.81"
$ws.Range("F8").Value = "'1.27"
$ws.Range("G8").Value = "'0.1010"
$ws.Range("H8").Value = 5
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "'009231"
$ws.Range("C9").Value = "鹏华安和混合C"
$ws.Range("D9").Value = "'5.33"
$ws.Range("E9").Value = "'34.45"
$ws.Range("F9").Value = "'1.67"
$ws.Range("G9").Value = "'0.0890"
$ws.Range("H9").Value = 2
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "'007854"
$ws.Range("C10").Value = "光大保德信景气先锋混合"
$ws.Range("D10").Value = "'1.82"
$ws.Range("E10").Value = "'72.20"
$ws.Range("F10").Value = "'4.11"
$ws.Range("G10").Value = "'0.0748"
$ws.Range("H10").Value = 10
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "'673141"
$ws.Range("C11").Value = "西部利得景程灵活配置混合A"
$ws.Range("D11").Value = "'2.09"
$ws.Range("E11").Value = "'86.50"
$ws.Range("F11").Value = "'3.33"
$ws.Range("G11").Value = "'0.0696"
$ws.Range("H11").Value = 6
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "'003165"
$ws.Range("C12").Value = "鹏华弘嘉灵活配置混合A"
$ws.Range("D12").Value = "'1.53"
$ws.Range("E12").Value = "'93.95"
$ws.Range("F12").Value = "'3.59"
$ws.Range("G12").Value = "'0.0549"
$ws.Range("H12").Value = 3
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "'009668"
$ws.Range("C13").Value = "鹏华安庆混合C"
$ws.Range("D13").Value = "'2.36"
$ws.Range("E13").Value = "'38.92"
$ws.Range("F13").Value = "'1.81"
$ws.Range("G13").Value = "'0.0427"
$ws.Range("H13").Value = 3
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "'011284"
$ws.Range("C14").Value = "中信保诚龙腾精选混合"
$ws.Range("D14").Value = "'1.22"
$ws.Range("E14").Value = "'75.38"
$ws.Range("F14").Value = "'2.62"
$ws.Range("G14").Value = "'0.0320"
$ws.Range("H14").Value = 10
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "'006209"
$ws.Range("C15").Value = "中信保诚新蓝筹灵活配置混合"
$ws.Range("D15").Value = "'1.16"
$ws.Range("E15").Value = "'77.03"
$ws.Range("F15").Value = "'2.66"
$ws.Range("G15").Value = "'0.0309"
$ws.Range("H15").Value = 10
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "'003166"
$ws.Range("C16").Value = "鹏华弘嘉灵活配置混合C"
$ws.Range("D16").Value = "'0.56"
$ws.Range("E16").Value = "'93.95"
$ws.Range("F16").Value = "'3.59"
$ws.Range("G16").Value = "'0.0201"
$ws.Range("H16").Value = 3
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "'673143"
$ws.Range("C17").Value = "西部利得景程灵活配置混合C"
$ws.Range("D17").Value = "'0.60"
$ws.Range("E17").Value = "'86.50"
$ws.Range("F17").Value = "'3.33"
$ws.Range("G17").Value = "'0.0200"
$ws.Range("H17").Value = 6
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "'005417"
$ws.Range("C18").Value = "鹏华尊惠18个月定期开放混合C"
$ws.Range("D18").Value = "'0.56"
$ws.Range("E18").Value = "'37.81"
$ws.Range("F18").Value = "'1.27"
$ws.Range("G18").Value = "'0.0071"
$ws.Range("H18").Value = 5
# -------------------------------------------------------------------
# Step 2: prepend a "2022-Q1" row to the "总计" (totals) summary sheet.
# -------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Push the existing data rows (2021-Q4 ... 2020-Q4) down by one row to
# make room for the new "2022-Q1" entry at the top of the table.
$zj.Rows.Item(2).Insert()

# The row-insert leaves the shifted-down cells in row 2 with an
# inherited/blank style; clear that so B2:D2 fall back to the sheet's
# plain (un-styled) look, matching every other data cell in columns
# B-D on this sheet.
$zj.Range("B2:D2").ClearFormats()

# Re-apply the bold/bordered/centred "index column" formatting (copied
# from the still-intact row below) onto the new A2 cell.
$zj.Range("A3").Copy()
$zj.Range("A2").PasteSpecial(-4122)

# New summary row for 2022-Q1.
$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 17
$zj.Range("D2").Value = 2.13

# Renumber the index column for the rows that shifted down.
$zj.Range("A3").Value = 1
$zj.Range("A4").Value = 2
$zj.Range("A5").Value = 3
$zj.Range("A6").Value = 4
$zj.Range("A7").Value = 5
